$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.955.80"
$ws.Range("E2").Value = "  +1.71%  "

$ws.Range("D3").Value = "1.647.72"
$ws.Range("E3").Value = "  +1.88%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.94"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.50%  "

$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.64"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.15%  "

$ws.Range("E9").Value = "  +1.79%  "

$ws.Range("E10").Value = "  +0.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0871"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.56%  "

$ws.Range("D12").Value = "1.880.79"
$ws.Range("E12").Value = "  +1.86%  "

$ws.Range("D13").Value = "1.647.99"
$ws.Range("E13").Value = "  +1.81%  "

$ws.Range("E14").Value = "  +1.55%  "

$ws.Range("E15").Value = "  +2.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.75"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.83%  "

$ws.Range("D17").Value = "27.971.49"
$ws.Range("E17").Value = "  +1.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.09"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.91%  "

$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.70"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +7.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.38"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.74"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.92"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.73"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.20%  "

$ws.Range("E28").Value = "  +0.28%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("E30").Value = "  +1.43%  "

$ws.Range("E31").Value = "  +0.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.20%  "

$ws.Range("D33").Value = "1.454.87"
$ws.Range("E33").Value = "  +1.05%  "

$ws.Range("E34").Value = "  +2.07%  "

$ws.Range("E35").Value = "  +2.06%  "

$ws.Range("E36").Value = "  -0.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.888"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.40%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0169"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.563"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.32%  "

$ws.Range("E40").Value = "  -2.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.44"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.45%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("E44").Value = "  +0.38%  "

$ws.Range("E45").Value = "  +0.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.38"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.24%  "

$ws.Range("E47").Value = "  +5.95%  "

$ws.Range("D48").Value = "1.789.50"
$ws.Range("E48").Value = "  +1.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.95"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.88%  "

$ws.Range("E50").Value = "  +1.50%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0999"
$ws.Range("E51").Value = "  -5.01%  "
